$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ADBE"
$ws.Range("B2").Value = 532.23
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 0.3615851006677071
$ws.Range("E2").Value = 0.8297872340425532
$ws.Range("F2").Value = 0.5447553259418354
$ws.Range("G2").Value = 0.9787234042553191
$ws.Range("H2").Value = 0.4050422386483632
$ws.Range("I2").Value = 0.9787234042553191
$ws.Range("J2").Value = 0.1723127753303966
$ws.Range("K2").Value = 0.9574468085106383
$ws.Range("L2").Value = 0.9361702127659575

$ws.Range("A3").Value = "AAL"
$ws.Range("B3").Value = 18.68
$ws.Range("C3").Value = 53
$ws.Range("D3").Value = 0.3248226950354609
$ws.Range("E3").Value = 0.8085106382978724
$ws.Range("F3").Value = 0.1099227569815806
$ws.Range("G3").Value = 0.6808510638297872
$ws.Range("H3").Value = 0.4424710424710425
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 0.1997430956968529
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.8723404255319149

$ws.Range("A4").Value = "ALGN"
$ws.Range("B4").Value = 360.62
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 0.3889231243259899
$ws.Range("E4").Value = 0.8936170212765957
$ws.Range("F4").Value = 0.5307751082434842
$ws.Range("G4").Value = 0.9574468085106383
$ws.Range("H4").Value = 0.05645232166398118
$ws.Range("I4").Value = 0.4468085106382979
$ws.Range("J4").Value = 0.1845355406648272
$ws.Range("K4").Value = 0.9787234042553191
$ws.Range("L4").Value = 0.8191489361702128

$ws.Range("A5").Value = "AMAT"
$ws.Range("B5").Value = 144.51
$ws.Range("C5").Value = 6
$ws.Range("D5").Value = 0.592396694214876
$ws.Range("E5").Value = 0.9787234042553191
$ws.Range("F5").Value = 0.3113430127041741
$ws.Range("G5").Value = 0.851063829787234
$ws.Range("H5").Value = 0.2735524808319378
$ws.Range("I5").Value = 0.9148936170212765
$ws.Range("J5").Value = 0.06163679106670572
$ws.Range("K5").Value = 0.4468085106382979
$ws.Range("L5").Value = 0.7978723404255319

$ws.Range("A6").Value = "AOS"
$ws.Range("B6").Value = 76.51
$ws.Range("C6").Value = 13
$ws.Range("D6").Value = 0.2835094782754572
$ws.Range("E6").Value = 0.6595744680851063
$ws.Range("F6").Value = 0.2460912052117266
$ws.Range("G6").Value = 0.8297872340425532
$ws.Range("H6").Value = 0.153127354935946
$ws.Range("I6").Value = 0.8085106382978724
$ws.Range("J6").Value = 0.1122256141881086
$ws.Range("K6").Value = 0.8297872340425532
$ws.Range("L6").Value = 0.7819148936170213

$ws.Range("A7").Value = "AAPL"
$ws.Range("B7").Value = 193.73
$ws.Range("C7").Value = 5
$ws.Range("D7").Value = 0.3237444482405192
$ws.Range("E7").Value = 0.7872340425531915
$ws.Range("F7").Value = 0.4521400194887939
$ws.Range("G7").Value = 0.9148936170212765
$ws.Range("H7").Value = 0.1701497946363855
$ws.Range("I7").Value = 0.851063829787234
$ws.Range("J7").Value = 0.07056808134394332
$ws.Range("K7").Value = 0.5319148936170213
$ws.Range("L7").Value = 0.7712765957446808

$ws.Range("A8").Value = "ABC"
$ws.Range("B8").Value = 192.39
$ws.Range("C8").Value = 5
$ws.Range("D8").Value = 0.3665033027913913
$ws.Range("E8").Value = 0.851063829787234
$ws.Range("F8").Value = 0.1667070952092176
$ws.Range("G8").Value = 0.7872340425531915
$ws.Range("H8").Value = 0.149420480344127
$ws.Range("I8").Value = 0.7659574468085105
$ws.Range("J8").Value = 0.07950847267422279
$ws.Range("K8").Value = 0.6808510638297872
$ws.Range("L8").Value = 0.7712765957446808

$ws.Range("A9").Value = "AMP"
$ws.Range("B9").Value = 352.46
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = 0.4792462332647835
$ws.Range("E9").Value = 0.9361702127659575
$ws.Range("F9").Value = 0.06428722408430709
$ws.Range("G9").Value = 0.5319148936170213
$ws.Range("H9").Value = 0.1430146581917238
$ws.Range("I9").Value = 0.7234042553191489
$ws.Range("J9").Value = 0.1212342929855255
$ws.Range("K9").Value = 0.8936170212765957
$ws.Range("L9").Value = 0.7712765957446808

$ws.Range("A10").Value = "ALLE"
$ws.Range("B10").Value = 125.6
$ws.Range("C10").Value = 7
$ws.Range("D10").Value = 0.2434412434412434
$ws.Range("E10").Value = 0.5957446808510638
$ws.Range("F10").Value = 0.106510439608845
$ws.Range("G10").Value = 0.6595744680851063
$ws.Range("H10").Value = 0.2296847464264733
$ws.Range("I10").Value = 0.8723404255319148
$ws.Range("J10").Value = 0.1149578339991122
$ws.Range("K10").Value = 0.8723404255319148
$ws.Range("L10").Value = 0.75

$ws.Range("A11").Value = "ALK"
$ws.Range("B11").Value = 53.42
$ws.Range("C11").Value = 18
$ws.Range("D11").Value = 0.3038808884549671
$ws.Range("E11").Value = 0.7446808510638298
$ws.Range("F11").Value = 0.08137651821862346
$ws.Range("G11").Value = 0.574468085106383
$ws.Range("H11").Value = 0.2605002359603585
$ws.Range("I11").Value = 0.8936170212765957
$ws.Range("J11").Value = 0.09804727646454281
$ws.Range("K11").Value = 0.7872340425531915
$ws.Range("L11").Value = 0.75
